# Add a new worksheet "Klay Thompson" with his stats, then wire the
# "final" summary sheet's row 3 to pull from it (mirrors the existing
# Stephen Curry / Draymond Green pattern).

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Klay Thompson" sheet right after "Draymond Green" ---
$draymond = $wb.Worksheets.Item("Draymond Green")
$klay = $wb.Worksheets.Add($null, $draymond)
$klay.Name = "Klay Thompson"

# --- 2. Header row (row 1) ---
$klay.Range("A1").Value = "TCA"
$klay.Range("B1").Value = "T3A"
$klay.Range("C1").Value = "TCI"
$klay.Range("D1").Value = "TLI"
$klay.Range("E1").Value = "Puntos"
$klay.Range("F1").Value = "%FG"
$klay.Range("G1").Value = "%eFG"
$klay.Range("H1").Value = "%TS"
$klay.Range("I1").Value = "Valoration"

# --- 3. Game stats (row 2) ---
$klay.Range("A2").Value = 51
$klay.Range("B2").Value = 6
$klay.Range("C2").Value = 256
$klay.Range("D2").Value = 45
$klay.Range("E2").Value = 153
$klay.Range("F2").Value = 19.92
$klay.Range("G2").Value = 21.09
$klay.Range("H2").Value = 27.74
$klay.Range("I2").Value = -77

# --- 4. Averages (row 3) ---
$klay.Range("A3").Formula = "=AVERAGE(A2:A2)"
$klay.Range("B3").Formula = "=AVERAGE(B2:B2)"
$klay.Range("C3").Formula = "=AVERAGE(C2:C2)"
$klay.Range("D3").Formula = "=AVERAGE(D2:D2)"
$klay.Range("E3").Formula = "=AVERAGE(E2:E2)"
$klay.Range("F3").Formula = "=AVERAGE(F2:F2)"
$klay.Range("G3").Formula = "=AVERAGE(G2:G2)"
$klay.Range("H3").Formula = "=AVERAGE(H2:H2)"
$klay.Range("I3").Formula = "=AVERAGE(I2:I2)"
$klay.Range("J3").Value = "promedios"

# --- 5. Wire the "final" sheet row 3 to the new Klay Thompson averages ---
$final = $wb.Worksheets.Item("final")
$final.Range("B3").Formula = "='Klay Thompson'!A3"
$final.Range("C3").Formula = "='Klay Thompson'!B3"
$final.Range("D3").Formula = "='Klay Thompson'!C3"
$final.Range("E3").Formula = "='Klay Thompson'!D3"
$final.Range("F3").Formula = "='Klay Thompson'!E3"
$final.Range("G3").Formula = "='Klay Thompson'!F3"
$final.Range("H3").Formula = "='Klay Thompson'!G3"
$final.Range("I3").Formula = "='Klay Thompson'!H3"
$final.Range("J3").Formula = "='Klay Thompson'!I3"

# Restore "final" as the active/selected sheet (it was the original tab
# shown before this edit; creating new sheets shifts the active tab).
$final.Activate() | Out-Null
$final.Range("B3").Select() | Out-Null

